# Add I0 and IF columns to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy style (bold font + border + alignment) from existing
# header cell (H1) to the new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-41
$data = @{
    2  = @(9, 9)
    3  = @(6, 6)
    4  = @(6, 6)
    5  = @(7, 7)
    6  = @(6, 7)
    7  = @(6, 7)
    8  = @(6, 7)
    9  = @(7, 7)
    10 = @(5, 6)
    11 = @(7, 8)
    12 = @(8, 9)
    13 = @(8, 8)
    14 = @(9, 9)
    15 = @(6, 7)
    16 = @(6, 7)
    17 = @(6, 7)
    18 = @(6, 7)
    19 = @(7, 7)
    20 = @(8, 8)
    21 = @(6, 7)
    22 = @(5, 6)
    23 = @(6, 6)
    24 = @(7, 8)
    25 = @(1, 2)
    26 = @(8, 8)
    27 = @(7, 7)
    28 = @(5, 6)
    29 = @(4, 5)
    30 = @(1, 4)
    31 = @(1, 4)
    32 = @(1, 5)
    33 = @(1, 5)
    34 = @(1, 4)
    35 = @(1, 3)
    36 = @(1, 3)
    37 = @(5, 6)
    38 = @(5, 6)
    39 = @(7, 7)
    40 = @(5, 6)
    41 = @(6, 6)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 9).Value = $values[0]
    $ws.Cells.Item($row, 10).Value = $values[1]
}
